# Fixing the big mistake: correct Total/Community/IGA consumption values
# in the monthly consumption table (rows 2-13 correspond to months 1-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 18069.38534810002;  D = 1401.702012583333;  E = 3133.176648116666 }
    3  = @{ B = 16907.53635013335;  D = 1314.817776483333;  E = 2813.500064166667 }
    4  = @{ B = 18226.97401928335;  D = 1413.251348733333;  E = 3000.74609245 }
    5  = @{ B = 17567.87675318335;  D = 1348.838251616667;  E = 3055.1659707 }
    6  = @{ B = 18159.73853466669;  D = 1388.288883433333;  E = 3132.18378535 }
    7  = @{ B = 17709.04371958335;  D = 1345.960124616667;  E = 2847.8812505 }
    8  = @{ B = 18113.95061670002;  D = 1394.813221483333;  E = 3077.713966683334 }
    9  = @{ B = 18167.91664681669;  D = 1400.696152216667;  E = 2903.36318005 }
    10 = @{ B = 17379.53084181668;  D = 1323.102249716667;  E = 3020.723615683334 }
    11 = @{ B = 18184.11341070002;  D = 1402.298082233333;  E = 3118.2647175 }
    12 = @{ B = 17567.74403295002;  D = 1347.85230415;      E = 3019.543172833333 }
    13 = @{ B = 17639.10509651668;  D = 1368.642479616667;  E = 2902.3091868 }
}

foreach ($row in $data.Keys) {
    $ws.Range("B$row").Value = $data[$row].B
    $ws.Range("D$row").Value = $data[$row].D
    $ws.Range("E$row").Value = $data[$row].E
}
